$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.205.50"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.706.77"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.27"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.47%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").Value = "2.707.14"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.29"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "3.220.04"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000187"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "68.144.35"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "2.723.87"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "369.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.58"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.48"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.08"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").Value = "2.854.90"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "574.85"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.14"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.33%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.59"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.376"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.87"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.38"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.98"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.87"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").Value = "0.0₆0308"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.592"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "154.14"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.89"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0785"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.04%  "
